$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: replace "Unaligned Augmented" placeholder row with real data for the
# unaligned model (style carried over from the "Aligned/Unaligned" comparison block above).
$ws.Range("A34").Value = "Unaligned_CL9_DL1_nobias_200Epoch_0.0001LR5Batch1keep0WD"
$ws.Range("B34").Value = 0.60563551777400004
$ws.Range("C34").Value = 0.28565417987199998
$ws.Range("A34:C34").Interior.Pattern = -4142  # xlPatternNone

# Row 35: "Unaligned Unaugmented" label row now highlighted with the yellow fill style.
$ws.Range("A35").Value = "Unaligned Unaugmented"
$ws.Range("B35").Value = "?"
$ws.Range("C35").Value = "?"
$ws.Range("A35:C35").Interior.Color = 65535

# Row 36: replace "Aligned Augmented" placeholder row with real data for the
# aligned model.
$ws.Range("A36").Value = "Aligned_CL9_DL1_nobias_200Epoch_0.0001LR5Batch1keep0WD"
$ws.Range("B36").Value = 0.90420036555600003
$ws.Range("C36").Value = 0.089764793522200001
$ws.Range("A36:C36").Interior.Pattern = -4142  # xlPatternNone

# Row 37: "Aligned Unaugmented" label row now highlighted with the yellow fill style.
$ws.Range("A37").Value = "Aligned Unaugmented"
$ws.Range("B37").Value = "?"
$ws.Range("C37").Value = "?"
$ws.Range("A37:C37").Interior.Color = 65535

# Update the view: move the selection to B33 (also clears the old A25 scroll anchor).
$ws.Range("B33").Select()
